$wb = $excel.ActiveWorkbook

# --- CaseDetailStat (sheet 5): replace the 4-column file/sample/case/study
# count table with a 6-column file-association table -----------------------
$wsCase = $wb.Worksheets.Item("CaseDetailStat")

# Clear the old table and drop in the new headers + data row.
$wsCase.Cells.Clear() | Out-Null

$wsCase.Range("A1").Value = "File Name"
$wsCase.Range("B1").Value = "File Type"
$wsCase.Range("C1").Value = "Association"
$wsCase.Range("D1").Value = "Description"
$wsCase.Range("E1").Value = "Format"
$wsCase.Range("F1").Value = "Size"

$wsCase.Range("A2").Value = "CCB010072.pdf"
$wsCase.Range("B2").Value = "Pathology Report"
$wsCase.Range("C2").Value = "diagnosis"
$wsCase.Range("D2").Value = ""
$wsCase.Range("E2").Value = "pdf"
# Size is stored as text (not a number) in the source data. Prefix with an
# apostrophe so Excel keeps it as text instead of auto-detecting a number.
$wsCase.Range("F2").Value = "'57.732421875"

$wsCase.Range("A1:F2").Columns.AutoFit() | Out-Null

# --- CaseDetailStat_Message (sheet 6): the 3rd logged query (row 28) now
# reflects the new file-association cypher query instead of the old
# file/sample/case/study count query ----------------------------------------
$wsCaseMsg = $wb.Worksheets.Item("CaseDetailStat_Message")
$wsCaseMsg.Range("A28").Value = "MATCH (f:file)-[*]->(c:case) WITH DISTINCT(f) AS f, c MATCH (f)-->(parent) WHERE c.case_id IN ['NCATS-COP01CCB010072'] RETURN f.file_name AS ``File Name`` ,f.file_type AS ``File Type``,head(labels(parent)) AS ``Association``, f.file_description AS ``Description``,f.file_format AS Format,((f.file_size)/1024) AS Size"

# --- view state: CaseDetailStat becomes the active/selected sheet ----------
$wsCase.Activate()
$wsCase.Range("C9").Select() | Out-Null

$ws1 = $wb.Worksheets.Item("CypherOutput")
$ws1.Activate()
$ws1.Range("H8").Select() | Out-Null

$ws4 = $wb.Worksheets.Item("StatOutput_Message")
$ws4.Activate()
$ws4.Range("A8").Select() | Out-Null

$wsCaseMsg.Activate()
$wsCaseMsg.Range("A28").Select() | Out-Null

$wsCase.Activate()
